$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the diff.
# Leading apostrophe forces Excel to store numeric-looking strings as text
# (preserving formats like "1.010" or "326.03" exactly as authored),
# matching the original inlineStr (text) cell type used throughout column D.
$ws.Range('D2').Value = '29.510.77'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.906.90'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('D4').Value = "'1.005"
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = "'326.03"
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('D7').Value = "'0.4839"
$ws.Range('E7').Value = '  +3.38%  '
$ws.Range('D8').Value = "'0.4075"
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').Value = "'0.08144"
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').Value = "'1.010"
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('D11').Value = "'23.46"
$ws.Range('E11').Value = '  +5.13%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.902.80'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'6.025"
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').Value = "'7.097"
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = "'90.30"
$ws.Range('E15').Value = '  +1.42%  '
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = "'0.06765"
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('D18').Value = "'0.00001045"
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').Value = "'17.71"
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = '29.534.14'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = "'5.614"
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').Value = "'11.84"
$ws.Range('E23').Value = '  +2.55%  '
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('D25').Value = '2.137.18'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = "'154.43"
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').Value = "'6.312"
$ws.Range('E28').Value = '  +10.72%  '
$ws.Range('D29').Value = "'2.107"
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('D30').Value = "'119.10"
$ws.Range('E30').Value = '  +1.87%  '
$ws.Range('E31').Value = '  -3.31%  '
$ws.Range('D32').Value = "'0.09563"
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = "'5.530"
$ws.Range('E33').Value = '  +2.86%  '
$ws.Range('D34').Value = "'1.396"
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('D36').Value = "'0.02271"
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').Value = "'0.06112"
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').Value = "'1.173"
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').Value = "'0.5956"
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('D40').Value = "'10.74"
$ws.Range('E40').Value = '  +6.26%  '
$ws.Range('D41').Value = "'7.938"
$ws.Range('D42').Value = "'0.1856"
$ws.Range('E42').Value = '  +1.11%  '
$ws.Range('D43').Value = "'2.442"
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'1.278"
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.07711"
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').Value = "'0.5573"
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('D48').Value = "'1.956"
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').Value = "'115.15"
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').Value = "'72.60"
$ws.Range('E50').Value = '  +1.91%  '
$ws.Range('D51').Value = "'1.053"
$ws.Range('E51').Value = '  +2.40%  '
